$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need Text format so Excel
# stores the literal digits (e.g. "556.70") instead of coercing to a float
# (which would normalize to 556.7 / lose trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '72.712.12'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '4.048.99'
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '556.70'
$ws.Range("E5").Value = '  +3.40%  '
$ws.Range("D6").Value = '152.50'
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("D7").Value = '4.043.61'
$ws.Range("E7").Value = '  +0.69%  '
$ws.Range("D8").Value = '0.696'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D10").Value = '0.760'
$ws.Range("E10").Value = '  +1.14%  '
$ws.Range("D11").Value = '0.173'
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("D12").Value = '53.93'
$ws.Range("E12").Value = '  +12.57%  '
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("D14").Value = '11.02'
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("D15").Value = '4.692.89'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("D16").Value = '4.056.72'
$ws.Range("E16").Value = '  +1.14%  '
$ws.Range("D17").Value = '14.50'
$ws.Range("E17").Value = '  +2.60%  '
$ws.Range("D18").Value = '20.80'
$ws.Range("E18").Value = '  +1.21%  '
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").Value = '72.690.38'
$ws.Range("E21").Value = '  +0.93%  '
$ws.Range("D22").Value = '449.39'
$ws.Range("E22").Value = '  +3.98%  '
$ws.Range("D23").Value = '98.05'
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("E24").Value = '  -0.84%  '
$ws.Range("D25").Value = '4.36'
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("D26").Value = '14.77'
$ws.Range("E26").Value = '  +1.52%  '
$ws.Range("D27").Value = '4.23'
$ws.Range("E27").Value = '  +11.89%  '
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("D29").Value = '10.93'
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("D30").Value = '5.95'
$ws.Range("E30").Value = '  +1.59%  '
$ws.Range("D31").Value = '37.41'
$ws.Range("E31").Value = '  +1.11%  '
$ws.Range("D32").Value = '7.87'
$ws.Range("E32").Value = '  +15.47%  '
$ws.Range("E33").Value = '  +3.22%  '
$ws.Range("D34").Value = '13.71'
$ws.Range("E34").Value = '  +1.85%  '
$ws.Range("D35").Value = '694.47'
$ws.Range("E35").Value = '  +2.14%  '
$ws.Range("D36").Value = '48.79'
$ws.Range("E36").Value = '  +15.17%  '
$ws.Range("D37").Value = '67.39'
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("D38").Value = '0.454'
$ws.Range("E38").Value = '  +6.32%  '
$ws.Range("D39").Value = '0.0₃0881'
$ws.Range("E39").Value = '  +6.33%  '
$ws.Range("E40").Value = '  -3.06%  '
$ws.Range("D41").Value = '3.45'
$ws.Range("E41").Value = '  -2.37%  '
$ws.Range("E42").Value = '  -2.01%  '
$ws.Range("D43").Value = '11.25'
$ws.Range("E43").Value = '  +17.88%  '
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("E48").Value = '  +3.13%  '
$ws.Range("E49").Value = '  +7.83%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '3.12'
$ws.Range("E50").Value = '  +2.94%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = '3.37'
$ws.Range("E51").Value = '  -0.05%  '
